# Updated symbol list on Sat Dec 17 08:43:05 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    # Force a Text number format so numeric-looking strings (with leading/trailing
    # zeros, exact precision, etc.) are preserved exactly as authored, matching
    # the original inline-string cell contents.
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "235.37"
Set-TextValue "D3"  "22.33"
Set-TextValue "D4"  "5.421"
Set-TextValue "D7"  "6.478"
Set-TextValue "D8"  "1.069"
Set-TextValue "D9"  "0.7833"
Set-TextValue "D10" "0.1396"
Set-TextValue "D11" "0.07400"
Set-TextValue "D12" "0.03186"
Set-TextValue "D13" "0.02956"
Set-TextValue "D14" "0.09258"
Set-TextValue "D15" "0.001665"
Set-TextValue "D16" "3.256"
Set-TextValue "D17" "0.04765"
Set-TextValue "D19" "0.006215"
Set-TextValue "D20" "0.005113"
Set-TextValue "D21" "0.001051"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.903"
Set-TextValue "D24" "2.135"
Set-TextValue "D27" "0.0004991"
Set-TextValue "D40" "0.04044"
Set-TextValue "D41" "0.006989"
Set-TextValue "D42" "0.1042"
Set-TextValue "D43" "0.002701"
Set-TextValue "D44" "0.009265"
Set-TextValue "D48" "0.03970"

# --- Volume(1h) label (column E) text tweaks ---
$ws.Range("E23").Value = "22LEOLEOBestin24h"
$ws.Range("E27").Value = "26UpBotsUBXT"

# --- Rows 42/43 swapped: CEJI and BKEXToken traded places ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
